$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-59 (A:rank, B:business_status, C:name, D:rating, E:user_ratings_total)
$data = @(
    ,@(44, "OPERATIONAL", "10 East Arts HUB", 5, 8)
    ,@(3, "OPERATIONAL", "Art On Main Gallery And Gifts", 5, 11)
    ,@(53, "OPERATIONAL", "Art To Remember", 3.9, 7)
    ,@(15, "OPERATIONAL", "ArtMix", 5, 8)
    ,@(0, "OPERATIONAL", "Arts Council of Indianapolis", 4.8, 12)
    ,@(6, "OPERATIONAL", "Arts for Lawrence", 4.9, 19)
    ,@(14, "OPERATIONAL", "Arts for Learning", 0, 0)
    ,@(32, "OPERATIONAL", "Carmel International Arts Festival", 4.8, 18)
    ,@(21, "OPERATIONAL", "Cat Head Press: Printshop and Artist Cooperative", 5, 26)
    ,@(9, "OPERATIONAL", "Christel De Haan Fine Arts Center", 4.6, 15)
    ,@(39, "OPERATIONAL", "Clowes Memorial Hall", 4.7, 797)
    ,@(50, "OPERATIONAL", "Creative Energy Arts", 0, 0)
    ,@(17, "OPERATIONAL", "Creative Grounds Fine Arts Academy", 5, 15)
    ,@(46, "OPERATIONAL", "Cultural Arts Gallery", 5, 1)
    ,@(29, "OPERATIONAL", "Dance Arts", 4.8, 10)
    ,@(54, "OPERATIONAL", "Dance Magic Performing Arts", 0, 0)
    ,@(38, "OPERATIONAL", "Edison School Of The Arts", 3.6, 36)
    ,@(8, "OPERATIONAL", "Gallery 924", 4.9, 15)
    ,@(12, "OPERATIONAL", "Garfield Park Arts Center (GPAC)", 4.6, 119)
    ,@(58, "OPERATIONAL", "Glass Arts Indiana, Inc.", 4.8, 16)
    ,@(25, "OPERATIONAL", "Harrison Center", 4.8, 134)
    ,@(40, "OPERATIONAL", "Herron School of Art and Design", 4.8, 20)
    ,@(30, "OPERATIONAL", "Herron School of Art and Design - Eskenazi Fine Arts Center", 4.5, 2)
    ,@(10, "OPERATIONAL", "High Frequency Arts", 5, 3)
    ,@(47, "OPERATIONAL", "Ignition Arts, LLC", 4.9, 8)
    ,@(16, "OPERATIONAL", "Incredible Art Department", 0, 0)
    ,@(2, "OPERATIONAL", "Indiana Arts Commission", 0, 0)
    ,@(34, "OPERATIONAL", "Indiana Fine Arts Academy", 4, 4)
    ,@(49, "OPERATIONAL", "Indiana Performing Arts Centre", 4.8, 6)
    ,@(5, "OPERATIONAL", "Indianapolis Art Center", 4.7, 112)
    ,@(56, "OPERATIONAL", "Indianapolis Arts Chorale", 0, 0)
    ,@(19, "OPERATIONAL", "Indianapolis Arts Services", 5, 1)
    ,@(18, "OPERATIONAL", "Indianapolis Artsgarden", 4.6, 146)
    ,@(51, "OPERATIONAL", "Jazz Arts Society of In Inc", 4, 1)
    ,@(20, "OPERATIONAL", "Korka International Arts", 0, 0)
    ,@(43, "OPERATIONAL", "Landmark Keystone Art Cinema", 4.4, 513)
    ,@(13, "OPERATIONAL", "Magdalena Gallery of Arts", 5, 1)
    ,@(36, "OPERATIONAL", "Mars Hill Arts Center", 4.9, 15)
    ,@(35, "OPERATIONAL", "Mooresville Arts Partnership", 4, 29)
    ,@(23, "OPERATIONAL", "Newfields", 4.7, 3896)
    ,@(31, "OPERATIONAL", "Palette Art Studio", 5, 7)
    ,@(48, "OPERATIONAL", "Performing Arts", 0, 0)
    ,@(55, "OPERATIONAL", "Red Barn Arts Collective", 5, 1)
    ,@(28, "OPERATIONAL", "Schrott Center for the Arts", 4.8, 106)
    ,@(1, "OPERATIONAL", "Sho Arts", 5, 1)
    ,@(59, "OPERATIONAL", "Sidney & Lois Eskenazi Fine Arts Center", 5, 1)
    ,@(7, "OPERATIONAL", "Southside Art League", 4.9, 14)
    ,@(37, "OPERATIONAL", "Steve Haigh Fine Art", 5, 1)
    ,@(27, "OPERATIONAL", "Studio Alchemy LLC", 5, 5)
    ,@(24, "OPERATIONAL", "The Art Studio of Carmel", 0, 0)
    ,@(11, "OPERATIONAL", "The Carmel Arts Council Children’s Art Gallery", 5, 1)
    ,@(26, "OPERATIONAL", "The Center for the Performing Arts", 4.7, 1193)
    ,@(4, "OPERATIONAL", "The Murphy Art Center", 5, 1)
    ,@(52, "OPERATIONAL", "The Palladium at the Center for the Performing Arts", 4.9, 197)
    ,@(33, "OPERATIONAL", "United Art & Education", 4.5, 132)
    ,@(41, "OPERATIONAL", "Warren Performing Arts Center", 4.5, 374)
    ,@(42, "OPERATIONAL", "carmel academy of the arts", 0, 0)
    ,@(57, "OPERATIONAL", "iDREAM Music And Arts Center", 0, 0)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# The refreshed data has one fewer row than before (60 -> 59); drop the now-stale last row
$ws.Rows.Item(60).Delete()

